$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 135
$ws.Range("I2").Value = 360
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 427
$ws.Range("N2").Value = 268
$ws.Range("P2").Value = 7
$ws.Range("R2").Value = 26
$ws.Range("S2").Value = 166
$ws.Range("T2").Value = 282
$ws.Range("U2").Value = 13
$ws.Range("V2").Value = 2312
$ws.Range("X2").Value = 2361
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 40
$ws.Range("AA2").Value = 12
